$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.558.59'
$ws.Range("E2").Value = '  -1.20%  '

$ws.Range("D3").Value = '2.262.68'
$ws.Range("E3").Value = '  +0.31%  '

$ws.Range("E4").Value = '  -0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.25'
$ws.Range("E5").Value = '  +0.70%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.641'
$ws.Range("E6").Value = '  +1.59%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '64.40'
$ws.Range("E7").Value = '  +1.92%  '

$ws.Range("E8").Value = '  -0.11%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.438'
$ws.Range("E9").Value = '  -0.64%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0955'
$ws.Range("E10").Value = '  -6.87%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '56.94'
$ws.Range("E11").Value = '  -0.56%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '26.41'
$ws.Range("E12").Value = '  +2.35%  '

$ws.Range("E13").Value = '  -1.26%  '

$ws.Range("D14").Value = '2.592.75'
$ws.Range("E14").Value = '  -0.07%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.95'
$ws.Range("E15").Value = '  -4.03%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.06'
$ws.Range("E16").Value = '  -1.72%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.822'
$ws.Range("E17").Value = '  -1.72%  '

$ws.Range("D18").Value = '2.261.46'
$ws.Range("E18").Value = '  -0.50%  '

$ws.Range("D19").Value = '43.394.41'
$ws.Range("E19").Value = '  -1.25%  '

$ws.Range("D20").Value = '0.0₃0967'
$ws.Range("E20").Value = '  -4.10%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.90'
$ws.Range("E21").Value = '  -0.15%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.10'
$ws.Range("E22").Value = '  +1.72%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '247.19'
$ws.Range("E23").Value = '  -1.29%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.81'
$ws.Range("E24").Value = '  +18.18%  '

$ws.Range("E25").Value = '  -0.01%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.42'
$ws.Range("E26").Value = '  +0.05%  '

$ws.Range("E27").Value = '  -2.21%  '

$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.70'
$ws.Range("E28").Value = '  -2.71%  '

$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '173.63'
$ws.Range("E29").Value = '  +1.10%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '21.68'
$ws.Range("E30").Value = '  +4.81%  '

$ws.Range("E31").Value = '  +3.60%  '

$ws.Range("E32").Value = '  -3.93%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.125'
$ws.Range("E33").Value = '  +1.08%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.94'
$ws.Range("E34").Value = '  +5.12%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0677'
$ws.Range("E35").Value = '  -0.15%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.91'
$ws.Range("E36").Value = '  +1.83%  '

$ws.Range("B37").Value = 'THORChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.40'
$ws.Range("E37").Value = '  -3.21%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.60'
$ws.Range("E38").Value = '  -4.96%  '

$ws.Range("E39").Value = '  -1.20%  '

$ws.Range("E40").Value = '  -2.61%  '

$ws.Range("E41").Value = '  -0.32%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.78'
$ws.Range("E42").Value = '  +6.95%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.49'
$ws.Range("E43").Value = '  +3.80%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.08'
$ws.Range("E44").Value = '  -1.21%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '96.62'
$ws.Range("E45").Value = '  -0.42%  '

$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.18'
$ws.Range("E46").Value = '  -0.59%  '

$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0939'
$ws.Range("E47").Value = '  -2.37%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.11'
$ws.Range("E48").Value = '  +6.40%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.000206'
$ws.Range("E49").Value = '  -1.53%  '

$ws.Range("D50").Value = '1.426.54'
$ws.Range("E50").Value = '  -0.52%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.26'
$ws.Range("E51").Value = '  -0.08%  '
